# Update quantity (F) and stock value (G) figures for several line items in the
# Companywise Stock Report, plus the associated "Sub Total:" / "Grand Total:"
# rows (column B) that aggregate them. Also corrects two swapped item-code
# values (B165 / B166).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F38").Value = 62
$ws.Range("G38").Value = 1936.26
$ws.Range("F43").Value = 18
$ws.Range("G43").Value = 3472.02
$ws.Range("F52").Value = 68
$ws.Range("G52").Value = 6360.72
$ws.Range("B63").Value = 38884.26
$ws.Range("F96").Value = 63
$ws.Range("G96").Value = 4430.16
$ws.Range("F102").Value = 92
$ws.Range("G102").Value = 7266.16
$ws.Range("F105").Value = 216
$ws.Range("G105").Value = 3572.64
$ws.Range("F119").Value = 8
$ws.Range("G119").Value = 906.08
$ws.Range("F121").Value = 129
$ws.Range("G121").Value = 2510.34
$ws.Range("B122").Value = 243747.55
$ws.Range("F136").Value = 44
$ws.Range("G136").Value = 4582.6
$ws.Range("B140").Value = 11498.57
$ws.Range("F155").Value = 115
$ws.Range("G155").Value = 2263.2
$ws.Range("F158").Value = 28
$ws.Range("G158").Value = 896.5599999999999
$ws.Range("B160").Value = 22059.09
$ws.Range("B165").Value = 53925
$ws.Range("B166").Value = 57756
$ws.Range("F235").Value = 27
$ws.Range("G235").Value = 2290.41
$ws.Range("B238").Value = 8928.24
$ws.Range("F333").Value = 149
$ws.Range("G333").Value = 16682.04
$ws.Range("F344").Value = 123
$ws.Range("G344").Value = 12435.3
$ws.Range("F346").Value = 48
$ws.Range("G346").Value = 9532.32
$ws.Range("F351").Value = 78
$ws.Range("G351").Value = 4612.14
$ws.Range("F360").Value = 50
$ws.Range("G360").Value = 9988.5
$ws.Range("F362").Value = 201
$ws.Range("G362").Value = 14132.31
$ws.Range("B370").Value = 347186.43
$ws.Range("F373").Value = 17
$ws.Range("G373").Value = 7593.05
$ws.Range("F378").Value = 0
$ws.Range("G378").Value = 0
$ws.Range("B379").Value = 27727.16
$ws.Range("F412").Value = 35
$ws.Range("G412").Value = 33943.35
$ws.Range("B413").Value = 33943.35
$ws.Range("F425").Value = 60
$ws.Range("G425").Value = 1475.4
$ws.Range("B428").Value = 46690
$ws.Range("F434").Value = 37
$ws.Range("G434").Value = 1499.98
$ws.Range("F437").Value = 51
$ws.Range("G437").Value = 1908.42
$ws.Range("F438").Value = 59
$ws.Range("G438").Value = 10999.96
$ws.Range("F440").Value = 59
$ws.Range("G440").Value = 3566.55
$ws.Range("F444").Value = 63
$ws.Range("G444").Value = 3446.1
$ws.Range("B445").Value = 42473.3
$ws.Range("F456").Value = 1
$ws.Range("G456").Value = 2456.93
$ws.Range("F458").Value = 1
$ws.Range("G458").Value = 2549.35
$ws.Range("B466").Value = 95314.21000000001
$ws.Range("F468").Value = 609
$ws.Range("G468").Value = 8191.05
$ws.Range("F473").Value = 336
$ws.Range("G473").Value = 5520.48
$ws.Range("F474").Value = 306
$ws.Range("G474").Value = 3919.86
$ws.Range("F475").Value = 376
$ws.Range("G475").Value = 7418.48
$ws.Range("F476").Value = 429
$ws.Range("G476").Value = 2822.82
$ws.Range("F478").Value = 118
$ws.Range("G478").Value = 2296.28
$ws.Range("F479").Value = 1001
$ws.Range("G479").Value = 6586.58
$ws.Range("F481").Value = 956
$ws.Range("G481").Value = 6204.44
$ws.Range("F482").Value = 436
$ws.Range("G482").Value = 5733.4
$ws.Range("F483").Value = 347
$ws.Range("G483").Value = 9126.1
$ws.Range("B486").Value = 104875.95
$ws.Range("F529").Value = 623
$ws.Range("G529").Value = 4236.4
$ws.Range("F530").Value = 363
$ws.Range("G530").Value = 2486.55
$ws.Range("F534").Value = 313
$ws.Range("G534").Value = 5173.89
$ws.Range("B537").Value = 40250.71
$ws.Range("F592").Value = 44
$ws.Range("G592").Value = 4511.76
$ws.Range("F596").Value = 7
$ws.Range("G596").Value = 195.51
$ws.Range("F599").Value = 115
$ws.Range("G599").Value = 3063.6
$ws.Range("B604").Value = 37609.58
$ws.Range("F606").Value = 115
$ws.Range("G606").Value = 15013.25
$ws.Range("F610").Value = 75
$ws.Range("G610").Value = 2040
$ws.Range("F612").Value = 15
$ws.Range("G612").Value = 408
$ws.Range("B613").Value = 62726.39
$ws.Range("F635").Value = 171
$ws.Range("G635").Value = 7383.78
$ws.Range("F637").Value = 85
$ws.Range("G637").Value = 3670.3
$ws.Range("B641").Value = 29203.74
$ws.Range("F707").Value = 59
$ws.Range("G707").Value = 6580.86
$ws.Range("B721").Value = 481199.11
$ws.Range("F732").Value = 42
$ws.Range("G732").Value = 4322.22
$ws.Range("B739").Value = 20017
$ws.Range("F746").Value = 89
$ws.Range("G746").Value = 3328.6
$ws.Range("B747").Value = 7929.38
$ws.Range("F792").Value = 850
$ws.Range("G792").Value = 138643.5
$ws.Range("F796").Value = 78
$ws.Range("G796").Value = 5265
$ws.Range("B797").Value = 165623.74
$ws.Range("B803").Value = 3058515.19
$ws.Range("B804").Value = 3058515.19
